$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (weekly refresh: a new data point pushed in,
# existing rows shift down one; the former last row's data ends up split
# across the new rows 38/39).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with this week's new record.
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44847
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100112045
$ws.Range("G6").Value = "Zapallo"
$ws.Range("H6").Value = "Camote"
$ws.Range("I6").Value = "1a nueva(o)"
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 870
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = 885
$ws.Range("N6").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O6").Value = "Perú"
$ws.Range("P6").Value = 885
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"
